# Fixed bottom row stacking 145s with 289s.
# A handful of rows on the Sub-Bundle_Data sheet represent SKUs that are not
# actually stackable bundle items (clips, drill bits, kits, etc). Their
# stacking-calculation columns (Qty, Dimensions, Unit) were left populated
# from a bad calc, so clear them out and zero the resulting Sub-Bundle Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sub-Bundle_Data")

$rows = @(206, 207, 208, 209, 210, 261)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Clear()           # C - Qty per pallet
    $ws.Cells.Item($r, 4).Value = 0         # D - Sub-Bundle Value
    $ws.Cells.Item($r, 5).Clear()           # E - dimension
    $ws.Cells.Item($r, 6).Clear()           # F - dimension
    $ws.Cells.Item($r, 7).Clear()           # G - dimension
    $ws.Cells.Item($r, 8).Clear()           # H - unit (UOM string)
}
